$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 311
$ws.Range("I4").Value = 397.33334
$ws.Range("J4").Value = 52
$ws.Range("K4").Value = 397.33334
$ws.Range("L4").Value = 52
$ws.Range("M4").Value = -283.33334
$ws.Range("N4").Value = -280

$ws.Range("H135").Value = 1007.3226
$ws.Range("I135").Value = 375.80435
$ws.Range("J135").Value = 2822.9375
$ws.Range("K135").Value = 3382.23915
$ws.Range("L135").Value = 25406.4375
$ws.Range("M135").Value = -847.2391499999999
$ws.Range("N135").Value = -30476.4375

$ws.Range("H137").Value = 3139.849
$ws.Range("I137").Value = 2605.524
$ws.Range("K137").Value = 7816.572
$ws.Range("M137").Value = -5266.572

$ws.Range("H140").Value = 47168
$ws.Range("J140").Value = 47168
$ws.Range("L140").Value = 47168
$ws.Range("N140").Value = -57528

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3855.979
$ws.Range("I32").Value = 2792.762
$ws.Range("J32").Value = 11975.091
$ws.Range("K32").Value = 2792.762
$ws.Range("L32").Value = 11975.091
$ws.Range("M32").Value = -2505.762
$ws.Range("N32").Value = -12549.091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 174342
$ws.Range("J42").Value = 174342
$ws.Range("L42").Value = 174342
$ws.Range("N42").Value = -174998

$ws.Range("H43").Value = 243684
$ws.Range("J43").Value = 243684
$ws.Range("L43").Value = 243684
$ws.Range("N43").Value = -244046

$ws.Range("H134").Value = 668.85455
$ws.Range("I134").Value = 585.2826
$ws.Range("J134").Value = 1096
$ws.Range("K134").Value = 1755.8478
$ws.Range("L134").Value = 3288
$ws.Range("M134").Value = 779.1522
$ws.Range("N134").Value = -8358

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1370.579
$ws.Range("I99").Value = 1271.3572
$ws.Range("J99").Value = 1648.4
$ws.Range("K99").Value = 1271.3572
$ws.Range("L99").Value = 1648.4
$ws.Range("M99").Value = 226.6428000000001
$ws.Range("N99").Value = -4644.4

$ws.Range("H126").Value = 1370.579
$ws.Range("I126").Value = 1271.3572
$ws.Range("J126").Value = 1648.4
$ws.Range("K126").Value = 3814.0716
$ws.Range("L126").Value = 4945.200000000001
$ws.Range("M126").Value = -1344.0716
$ws.Range("N126").Value = -9885.200000000001

$ws.Range("H132").Value = 1007.88574
$ws.Range("I132").Value = 841.03845
$ws.Range("J132").Value = 1489.8889
$ws.Range("K132").Value = 2523.11535
$ws.Range("L132").Value = 4469.6667
$ws.Range("M132").Value = 6.884649999999965
$ws.Range("N132").Value = -9529.6667

$ws.Range("H134").Value = 1695.6305
$ws.Range("I134").Value = 1412.5294
$ws.Range("J134").Value = 2497.75
$ws.Range("K134").Value = 4237.5882
$ws.Range("L134").Value = 7493.25
$ws.Range("M134").Value = -1702.5882
$ws.Range("N134").Value = -12563.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6802.3125
$ws.Range("I68").Value = 240
$ws.Range("J68").Value = 7739.7856
$ws.Range("K68").Value = 720
$ws.Range("L68").Value = 23219.3568
$ws.Range("M68").Value = 91
$ws.Range("N68").Value = -24841.3568

$ws.Range("H71").Value = 6802.3125
$ws.Range("I71").Value = 240
$ws.Range("J71").Value = 7739.7856
$ws.Range("K71").Value = 2160
$ws.Range("L71").Value = 69658.0704
$ws.Range("M71").Value = 1896
$ws.Range("N71").Value = -77770.0704

$ws.Range("H131").Value = 810.35803
$ws.Range("I131").Value = 419.9375
$ws.Range("J131").Value = 906.46155
$ws.Range("K131").Value = 1259.8125
$ws.Range("L131").Value = 2719.38465
$ws.Range("M131").Value = 3780.1875
$ws.Range("N131").Value = -12799.38465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 49
$ws.Range("I2").Value = 48.5
$ws.Range("J2").Value = 49.857143
$ws.Range("K2").Value = 48.5
$ws.Range("L2").Value = 49.857143
$ws.Range("M2").Value = 64.5
$ws.Range("N2").Value = -275.857143

$ws.Range("H102").Value = 1679.4615
$ws.Range("I102").Value = 1864.091
$ws.Range("J102").Value = 664
$ws.Range("K102").Value = 1864.091
$ws.Range("L102").Value = 664
$ws.Range("M102").Value = -242.0909999999999
$ws.Range("N102").Value = -3908

$ws.Range("H122").Value = 1941.8636
$ws.Range("I122").Value = 1973.3889
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5920.1667
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -3470.1667
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3653.7
$ws.Range("I40").Value = 3045.7273
$ws.Range("K40").Value = 3045.7273
$ws.Range("M40").Value = -2909.7273

$ws.Range("H46").Value = 474.83334
$ws.Range("I46").Value = 449.8
$ws.Range("K46").Value = 449.8
$ws.Range("M46").Value = -261.8

$ws.Range("H55").Value = 630.125
$ws.Range("I55").Value = 480
$ws.Range("J55").Value = 720.2
$ws.Range("K55").Value = 480
$ws.Range("L55").Value = 720.2
$ws.Range("M55").Value = -307
$ws.Range("N55").Value = -1066.2

$ws.Range("H122").Value = 4003.3572
$ws.Range("I122").Value = 3983.7
$ws.Range("J122").Value = 4052.5
$ws.Range("K122").Value = 11951.1
$ws.Range("L122").Value = 12157.5
$ws.Range("M122").Value = -9501.099999999999
$ws.Range("N122").Value = -17057.5

$ws.Range("H132").Value = 2724.2368
$ws.Range("I132").Value = 2286.1428
$ws.Range("J132").Value = 3950.9
$ws.Range("K132").Value = 6858.428400000001
$ws.Range("L132").Value = 11852.7
$ws.Range("M132").Value = -4328.428400000001
$ws.Range("N132").Value = -16912.7

$ws.Range("H133").Value = 23118
$ws.Range("J133").Value = 23118
$ws.Range("L133").Value = 23118
$ws.Range("N133").Value = -28178

$ws.Range("H134").Value = 30018.428
$ws.Range("J134").Value = 38025.8
$ws.Range("L134").Value = 38025.8
$ws.Range("N134").Value = -48165.8

$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140

$ws.Range("H136").Value = 2587.0454
$ws.Range("I136").Value = 1868.5714
$ws.Range("J136").Value = 3398.2258
$ws.Range("K136").Value = 5605.7142
$ws.Range("L136").Value = 10194.6774
$ws.Range("M136").Value = -3055.7142
$ws.Range("N136").Value = -15294.6774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1450
$ws.Range("I81").Value = 1450
$ws.Range("K81").Value = 2900
$ws.Range("M81").Value = -1839

$ws.Range("H84").Value = 1450
$ws.Range("I84").Value = 1450
$ws.Range("K84").Value = 14500
$ws.Range("M84").Value = -9196

$ws.Range("H132").Value = 736
$ws.Range("I132").Value = 610.9729599999999
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 1832.91888
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = 697.0811200000003
$ws.Range("N132").Value = -8810
